$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 (Alvaro)
    $ws.Range("B9").Value = 'Verificar'
    $ws.Range("D9").Value = 0
    $ws.Range("E9").Value = 0
    $ws.Range("F9").Value = 0
    $ws.Range("G9").Value = 0
    $ws.Range("H9").Value = 0

# Row 10 (Sotto ツ)
    $ws.Range("B10").Value = 'Verificar'
    $ws.Range("D10").Value = 3
    $ws.Range("E10").Value = 5
    $ws.Range("G10").Value = 12
    $ws.Range("H10").Value = 12

# Row 11 (Felipe)
    $ws.Range("B11").Value = 'Ok'
    $ws.Range("D11").Value = 16
    $ws.Range("E11").Value = 16
    $ws.Range("F11").Value = 16
    $ws.Range("G11").Value = 5

# Row 12 (Zenitsu愛)
    $ws.Range("B12").Value = 'Razoável'
    $ws.Range("D12").Value = 14
    $ws.Range("E12").Value = 16
    $ws.Range("F12").Value = 16
    $ws.Range("G12").Value = 16
    $ws.Range("H12").Value = 16

# Row 13 (danilo)
    $ws.Range("B13").Value = 'Verificar'
    $ws.Range("D13").Value = 0
    $ws.Range("E13").Value = 0
    $ws.Range("F13").Value = 0
    $ws.Range("G13").Value = 0

# Row 14 (PedrinhoR14)
    $ws.Range("B14").Value = 'Verificar'
    $ws.Range("D14").Value = 0
    $ws.Range("E14").Value = 12
    $ws.Range("F14").Value = 8

# Row 15 (Dockz)
    $ws.Range("E15").Value = 16
    $ws.Range("H15").Value = 16

# Row 16 (5C4RF4C3)
    $ws.Range("E16").Value = 1

# Row 17 (luba)
    $ws.Range("D17").Value = 9
    $ws.Range("E17").Value = 16
    $ws.Range("F17").Value = 16
    $ws.Range("H17").Value = 12

# Row 18 (ed)
    $ws.Range("B18").Value = 'Verificar'
    $ws.Range("D18").Value = 0
    $ws.Range("E18").Value = 15
    $ws.Range("F18").Value = 9
    $ws.Range("G18").Value = 11
    $ws.Range("H18").Value = 12

# Row 19 (caioba_)
    $ws.Range("B19").Value = 'Razoável'
    $ws.Range("D19").Value = 12
    $ws.Range("E19").Value = 11
    $ws.Range("F19").Value = 16
    $ws.Range("G19").Value = 14
    $ws.Range("H19").Value = 16

# Row 20 (BRS⚔️ASHURA)
    $ws.Range("D20").Value = 4
    $ws.Range("E20").Value = 0
    $ws.Range("F20").Value = 0
    $ws.Range("G20").Value = 0
    $ws.Range("H20").Value = 0

# Row 21 (john)
    $ws.Range("E21").Value = 12
    $ws.Range("F21").Value = 12
    $ws.Range("G21").Value = 12
    $ws.Range("H21").Value = 9

# Row 22 (Theus Carvalho)
    $ws.Range("B22").Value = 'Verificar'
    $ws.Range("D22").Value = 8
    $ws.Range("E22").Value = 16
    $ws.Range("F22").Value = 12
    $ws.Range("G22").Value = 13
    $ws.Range("H22").Value = 10

# Row 23 (mini.peka)
    $ws.Range("D23").Value = 0

# Row 24 (kauansin777)
    $ws.Range("E24").Value = 0
    $ws.Range("F24").Value = 0
    $ws.Range("G24").Value = 0
    $ws.Range("H24").Value = 0

# Row 25 (Luiz Fernando™)
    $ws.Range("B25").Value = 'Razoável'
    $ws.Range("D25").Value = 12
    $ws.Range("G25").Value = 12
    $ws.Range("H25").Value = 16

# Row 26 (dogmal)
    $ws.Range("B26").Value = 'Ok'
    $ws.Range("D26").Value = 16
    $ws.Range("E26").Value = 16
    $ws.Range("H26").Value = 16

# Row 27 (Rodolfos)
    $ws.Range("F27").Value = 8
    $ws.Range("G27").Value = 12

# Row 28 (polaris)
    $ws.Range("B28").Value = 'Razoável'
    $ws.Range("D28").Value = 14
    $ws.Range("E28").Value = 14
    $ws.Range("F28").Value = 9
    $ws.Range("G28").Value = 14

# Row 29 (Chetto)
    $ws.Range("B29").Value = 'Ok'
    $ws.Range("D29").Value = 16
    $ws.Range("E29").Value = 16
    $ws.Range("F29").Value = 8
    $ws.Range("G29").Value = 0
    $ws.Range("H29").Value = 0

# Row 30 (Teixeirazzqw)
    $ws.Range("D30").Value = 0

# Row 31 (WILLIAN)
    $ws.Range("B31").Value = 'Ok'
    $ws.Range("D31").Value = 16
    $ws.Range("E31").Value = 14
    $ws.Range("F31").Value = 16
    $ws.Range("G31").Value = 16
    $ws.Range("H31").Value = 16

# Row 32 (Mila)
    $ws.Range("B32").Value = 'Ok'
    $ws.Range("D32").Value = 16

# Row 33 (RaiNascimento)
    $ws.Range("B33").Value = 'Razoável'
    $ws.Range("D33").Value = 12
    $ws.Range("E33").Value = 14
    $ws.Range("F33").Value = 12
    $ws.Range("G33").Value = 13
    $ws.Range("H33").Value = 16

# Row 34 (tavin)
    $ws.Range("D34").Value = 14
    $ws.Range("E34").Value = 10
    $ws.Range("F34").Value = 16
    $ws.Range("G34").Value = 16

# Row 35 (^_^^_^)
    $ws.Range("B35").Value = 'Ok'
    $ws.Range("D35").Value = 16
    $ws.Range("F35").Value = 16
    $ws.Range("G35").Value = 16

# Row 36 (Bruno)
    $ws.Range("B36").Value = 'Razoável'
    $ws.Range("D36").Value = 12
    $ws.Range("E36").Value = 13
    $ws.Range("F36").Value = 11
    $ws.Range("G36").Value = 16
    $ws.Range("H36").Value = 15

# Row 37 (⭐O SENTINELA ⭐)
    $ws.Range("B37").Value = 'Verificar'
    $ws.Range("D37").Value = 0
    $ws.Range("E37").Value = 0
    $ws.Range("F37").Value = 0
    $ws.Range("G37").Value = 0
    $ws.Range("H37").Value = 0

# Row 38 (domador de but)
    $ws.Range("B38").Value = 'Razoável'
    $ws.Range("D38").Value = 13
    $ws.Range("E38").Value = 14
    $ws.Range("F38").Value = 15
    $ws.Range("G38").Value = 14
    $ws.Range("H38").Value = 14

# Row 39 (filho de duque)
    $ws.Range("D39").Value = 10

# Row 40 (O GUARDIÃO)
    $ws.Range("G40").Value = 15
    $ws.Range("H40").Value = 12

# Row 41 (luck)
    $ws.Range("B41").Value = 'Verificar'
    $ws.Range("D41").Value = 0
    $ws.Range("E41").Value = 0
    $ws.Range("F41").Value = 0
    $ws.Range("G41").Value = 0
    $ws.Range("H41").Value = 0

# Row 42 (Nico)
    $ws.Range("B42").Value = 'Ok'
    $ws.Range("D42").Value = 16
    $ws.Range("E42").Value = 16
    $ws.Range("F42").Value = 12
    $ws.Range("G42").Value = 16
    $ws.Range("H42").Value = 13

# Row 43 (Pedro PH)
    $ws.Range("B43").Value = 'Razoável'
    $ws.Range("D43").Value = 12
    $ws.Range("E43").Value = 8
    $ws.Range("F43").Value = 14
    $ws.Range("G43").Value = 12
    $ws.Range("H43").Value = 16

# Row 44 (DGJ-DAVI)
    $ws.Range("F44").Value = 16
    $ws.Range("G44").Value = 16

# Row 45 (Luciano)
    $ws.Range("E45").Value = 4
    $ws.Range("F45").Value = 0
    $ws.Range("G45").Value = 0
    $ws.Range("H45").Value = 0

# Row 47 (GabiMalvadeza)
    $ws.Range("B47").Value = 'Ok'
    $ws.Range("D47").Value = 16
    $ws.Range("E47").Value = 16
    $ws.Range("F47").Value = 16
    $ws.Range("G47").Value = 16
    $ws.Range("H47").Value = 16

# Row 49 (EDDIE)
    $ws.Range("B49").Value = 'Razoável'
    $ws.Range("D49").Value = 15

# Row 50 (WvCly)
    $ws.Range("B50").Value = 'Razoável'
    $ws.Range("D50").Value = 15
    $ws.Range("E50").Value = 8
    $ws.Range("F50").Value = 16
    $ws.Range("G50").Value = 13
    $ws.Range("H50").Value = 16

# Row 51 (nivelador)
    $ws.Range("B51").Value = 'Ok'
    $ws.Range("D51").Value = 16
    $ws.Range("E51").Value = 14
    $ws.Range("F51").Value = 16
    $ws.Range("G51").Value = 16
    $ws.Range("H51").Value = 15

# Row 52 (andrebts)
    $ws.Range("B52").Value = 'Verificar'
    $ws.Range("D52").Value = 8

# Row 53 (OneDePrata)
    $ws.Range("B53").Value = 'Ok'
    $ws.Range("D53").Value = 16
    $ws.Range("E53").Value = 16
    $ws.Range("G53").Value = 16

# Row 55 (gabiggoughost)
    $ws.Range("E55").Value = 16

# Row 56 (Grimmer 狼)
    $ws.Range("B56").Value = 'Verificar'
    $ws.Range("D56").Value = 8
    $ws.Range("E56").Value = 14
    $ws.Range("F56").Value = 7
    $ws.Range("G56").Value = 10
    $ws.Range("H56").Value = 6

# Row 57 (joão3:16)
    $ws.Range("F57").Value = 12
    $ws.Range("G57").Value = 15

# Row 59 (51 é pinga)
    $ws.Range("B59").Value = 'Razoável'
    $ws.Range("D59").Value = 12
    $ws.Range("E59").Value = 10
    $ws.Range("F59").Value = 15
    $ws.Range("G59").Value = 12
    $ws.Range("H59").Value = 0

# Row 60 (gabriel 3$)
    $ws.Range("F60").Value = 14
    $ws.Range("H60").Value = 16

# Row 61 (Asten Acady)
    $ws.Range("H61").Value = 4

# Row 63 (mathi❄️)
    $ws.Range("B63").Value = 'Ok'
    $ws.Range("D63").Value = 16
    $ws.Range("E63").Value = 16
